$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "82-55=27",
    "56-54=2",
    "41-6=35",
    "27-20=7",
    "28+47=75",
    "75+2=77",
    "47-6=41",
    "8+71=79",
    "46+15=61",
    "21-9=12",
    "19+58=77",
    "0+83=83",
    "42-10=32",
    "69-38=31",
    "26-2=24",
    "2+9=11",
    "19+64=83",
    "16+12=28",
    "37+1=38",
    "77+21=98",
    "14+17=31",
    "7+35=42",
    "73-43=30",
    "97-87=10",
    "7+89=96",
    "27+30=57",
    "84-79=5",
    "73-56=17",
    "67-15=52",
    "61-25=36",
    "4+66=70",
    "70-58=12",
    "28+44=72",
    "47+50=97",
    "97-63=34",
    "84-27=57",
    "22+30=52",
    "23+2=25",
    "52+6=58",
    "39+20=59",
    "10-0=10",
    "79-4=75",
    "72-17=55",
    "39+60=99",
    "75+23=98",
    "27-0=27",
    "61+31=92",
    "66+18=84",
    "32+49=81",
    "12+44=56",
    "84-67=17",
    "63-61=2",
    "75-11=64",
    "23+73=96",
    "97-30=67",
    "43-41=2",
    "50+49=99",
    "74-51=23",
    "41+30=71",
    "18+68=86",
    "1+75=76",
    "50+8=58",
    "31-11=20",
    "43+41=84",
    "78-26=52",
    "11+79=90",
    "90-4=86",
    "73+26=99",
    "88-55=33",
    "13+69=82",
    "6+46=52",
    "95-0=95",
    "79+19=98",
    "35+11=46",
    "60+29=89",
    "98-7=91",
    "24+28=52",
    "47-46=1",
    "86-76=10",
    "94-64=30",
    "64-5=59",
    "90-34=56",
    "43-16=27",
    "12-11=1",
    "20+46=66",
    "37+0=37",
    "28-4=24",
    "48+12=60",
    "25+73=98",
    "63-33=30",
    "33-30=3",
    "2+44=46",
    "48+29=77",
    "77+22=99",
    "37-12=25",
    "68-45=23",
    "71-56=15",
    "14+29=43",
    "50-16=34",
    "55+21=76"
)
$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        [void]$rng.MoveEnd(1, -1)
        $rng.Text = $values[$idx]
        $idx++
    }
}
Write-Output "Updated $idx cells"